# cs-en-us-022pct.xlsx weekly refresh: new crime data collected.
# Bumps the report "Volume/Number" and the covered week, and replaces the
# weekly crime-complaint figures (rows 16-27) with the newly collected
# counts / percentage changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helpers -----------------------------------------------------------
# Set a numeric cell's value while preserving the existing numeric style
# (re-applying the same NumberFormat the workbook already uses reuses the
# matching cellXf instead of minting a new one).
function Set-NumCell($ws, $addr, $value, $fmt) {
    $ws.Range($addr).Value = $value
    $ws.Range($addr).NumberFormat = $fmt
}

# Set a cell to a literal text value (e.g. the "0" / "***.*" placeholders)
# while adopting the exact style of a known-good template cell that
# already carries that style (copy/paste-special formats only).
function Set-TextCell($ws, $addr, $text, $styleSrc) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($styleSrc).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# --- header: volume number + week covered -------------------------------
$ws.Range("A8").Value = "Volume 30   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/5/2023  Through  6/11/2023"

# --- row 16 --------------------------------------------------------------
Set-NumCell  $ws "C16" 1    '#,##0'
Set-TextCell $ws "D16" "0"    "C14"
Set-TextCell $ws "E16" "***.*" "C14"
Set-NumCell  $ws "F16" 2    '#,##0'
Set-NumCell  $ws "G16" 1    '#,##0'
Set-NumCell  $ws "H16" 100  '#,##0.0;"-"#,##0.0'
Set-NumCell  $ws "I16" 6    '#,##0'
Set-NumCell  $ws "K16" 20   '#,##0.0;"-"#,##0.0'
Set-NumCell  $ws "L16" -14.285714285714 '#,##0.0;"-"#,##0.0'
Set-NumCell  $ws "M16" -33.333333333333 '#,##0.0;"-"#,##0.0'
Set-NumCell  $ws "N16" -89.830508474576 '#,##0.0;"-"#,##0.0'

# --- row 17 --------------------------------------------------------------
Set-NumCell  $ws "N17" -85  '#,##0.0;"-"#,##0.0'

# --- row 18 --------------------------------------------------------------
Set-TextCell $ws "C18" "0"  "C14"
Set-NumCell  $ws "M18" 50   '#,##0.0;"-"#,##0.0'

# --- row 19 --------------------------------------------------------------
Set-TextCell $ws "C19" "0"    "C14"
Set-TextCell $ws "D19" "0"    "C14"
Set-TextCell $ws "E19" "***.*" "C14"
Set-NumCell  $ws "F19" 4    '#,##0'
Set-NumCell  $ws "G19" 3    '#,##0'
Set-NumCell  $ws "H19" 33.333333333333 '#,##0.0;"-"#,##0.0'
Set-NumCell  $ws "M19" -48.148148148148 '#,##0.0;"-"#,##0.0'
Set-NumCell  $ws "N19" -78.125 '#,##0.0;"-"#,##0.0'

# --- row 21 (TOTAL, bold styles 17/18/19) --------------------------------
Set-NumCell  $ws "C21" 1    '#,##0'
Set-TextCell $ws "D21" "0"    "A21"
Set-TextCell $ws "E21" "***.*" "A21"
Set-NumCell  $ws "F21" 8    '#,##0'
Set-NumCell  $ws "G21" 6    '#,##0'
Set-NumCell  $ws "H21" 33.333333333333 '#,##0.00;"-"#,##0.00'
Set-NumCell  $ws "I21" 26   '#,##0'
Set-NumCell  $ws "K21" 23.809523809523 '#,##0.00;"-"#,##0.00'
Set-NumCell  $ws "L21" 4    '#,##0.00;"-"#,##0.00'
Set-NumCell  $ws "M21" -36.585365853658 '#,##0.00;"-"#,##0.00'
Set-NumCell  $ws "N21" -84.049079754601 '#,##0.00;"-"#,##0.00'

# --- row 24 --------------------------------------------------------------
Set-TextCell $ws "C24" "0"  "C14"
Set-NumCell  $ws "D24" 3    '#,##0'
Set-NumCell  $ws "E24" -100 '#,##0.0;"-"#,##0.0'
Set-NumCell  $ws "F24" 2    '#,##0'
Set-NumCell  $ws "G24" 6    '#,##0'
Set-NumCell  $ws "H24" -66.666666666666 '#,##0.0;"-"#,##0.0'
Set-NumCell  $ws "J24" 12   '#,##0'
Set-NumCell  $ws "K24" 16.666666666666 '#,##0.0;"-"#,##0.0'
Set-NumCell  $ws "M24" -50  '#,##0.0;"-"#,##0.0'

# --- row 25 --------------------------------------------------------------
Set-NumCell  $ws "C25" 1    '#,##0'
Set-NumCell  $ws "E25" 0    '#,##0.0;"-"#,##0.0'
Set-NumCell  $ws "G25" 2    '#,##0'
Set-NumCell  $ws "H25" 300  '#,##0.0;"-"#,##0.0'
Set-NumCell  $ws "I25" 17   '#,##0'
Set-NumCell  $ws "J25" 14   '#,##0'
Set-NumCell  $ws "K25" 21.428571428571 '#,##0.0;"-"#,##0.0'
Set-NumCell  $ws "L25" 112.5 '#,##0.0;"-"#,##0.0'
Set-NumCell  $ws "M25" 41.666666666666 '#,##0.0;"-"#,##0.0'

# --- row 27 --------------------------------------------------------------
Set-NumCell  $ws "C27" 1    '#,##0'
Set-NumCell  $ws "I27" 10   '#,##0'
Set-NumCell  $ws "K27" 233.333333333333 '#,##0.0;"-"#,##0.0'
Set-NumCell  $ws "L27" 900  '#,##0.0;"-"#,##0.0'
